# AfDD_2023_Annex_Table_Tab06.xlsx - "Add files via upload" update
# Re-applies the refreshed data values and fixes the mojibake'd
# accented characters in the Regional Economic Communities footnote.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tab06")

# --- Fix mis-encoded accented characters in the footnote text (A103) ---
$ws.Range("A103").Value = 'Regional Economic Communities:CEN-SAD = "Community of Sahel-Saharan States";COMESA = "Common Market for Eastern and Southern Africa";EAC = "East African Community";ECCAS = "Economic Community of Central African States";ECOWAS = "Economic Community of West African States";IGAD = "Intergovernmental Authority on Development";SADC = "Southern African Development Community";UMA = "Arab Maghreb Union";PALOP = "Países Africanos de Língua Oficial Portuguesa";ASEAN = "Association of Southeast Asian Nations";MERCOSUR = "Mercado Común del Sur".EU27 = "European Union (27 members)".OECD = "Organisation for Economic Co-operation and Development".'

# --- Refreshed data values (re-uploaded source recalculated these) ---
$ws.Range("D62").Value = 76.054260200000002
$ws.Range("E62").Value = 81.404715400000001
$ws.Range("H63").Value = 91.049872470588198
$ws.Range("I63").Value = 93.952597294117595
$ws.Range("C65").Value = 94.159400000000005
$ws.Range("G65").Value = 86.209238518518504
$ws.Range("I65").Value = 89.690868518518499
$ws.Range("C66").Value = 90.209619402985098
$ws.Range("E66").Value = 91.359260746268703
$ws.Range("H66").Value = 79.986195851851903
$ws.Range("D68").Value = 64.433781818181799
$ws.Range("E68").Value = 74.776758636363596
$ws.Range("C71").Value = 67.974553333333404
$ws.Range("E71").Value = 74.632891333333404
$ws.Range("E72").Value = 77.4485733333334
$ws.Range("E76").Value = 97.176111000000006
$ws.Range("I79").Value = 98.146425833333296
$ws.Range("F81").Value = 1.0029744444444399
$ws.Range("G81").Value = 94.470420555555606
$ws.Range("F82").Value = 0.92154190476191
$ws.Range("I82").Value = 74.578712619047593
$ws.Range("D83").Value = 96.736173333333397
$ws.Range("E83").Value = 97.182109393939399
$ws.Range("G83").Value = 91.962249999999997
$ws.Range("I86").Value = 79.9204890476191
$ws.Range("H87").Value = 82.605646538461599
$ws.Range("G89").Value = 95.951861944444403
$ws.Range("C90").Value = 99.108101428571402
$ws.Range("E90").Value = 98.965710000000101
$ws.Range("G90").Value = 97.481631818181896
$ws.Range("H90").Value = 97.137345909090897
$ws.Range("I90").Value = 97.744622727272699
$ws.Range("C91").Value = 70.356761612903199
$ws.Range("H91").Value = 49.7276545161291
$ws.Range("I93").Value = 83.904233333333394
$ws.Range("G94").Value = 91.108957333333393

# Rows 97-98 ("Africa, Fragile States" / "ROW, Fragile States") were
# recomputed from an updated underlying country sample.
$ws.Range("C97").Value = 73.668772692307698
$ws.Range("D97").Value = 69.921260769230798
$ws.Range("E97").Value = 77.852894230769195
$ws.Range("F97").Value = 0.88409269230769005
$ws.Range("G97").Value = 62.093014615384597
$ws.Range("H97").Value = 55.0898015384616
$ws.Range("I97").Value = 69.698148846153899
$ws.Range("J97").Value = 0.76230192307692002

$ws.Range("C98").Value = 93.047320769230794
$ws.Range("D98").Value = 92.761195384615405
$ws.Range("E98").Value = 93.315020769230799
$ws.Range("F98").Value = 0.99209307692308002
$ws.Range("G98").Value = 83.795082307692297
$ws.Range("H98").Value = 80.6872969230769
$ws.Range("I98").Value = 86.988493846153901
$ws.Range("J98").Value = 0.92114384615385003
